# Update the "working experience" table with the new entries described in the
# commit message: replace the old "Associate Professor / Universidad El
# Bosque" row with two new consultancy rows (RedPapaz, Protect Children).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: External Scientific Consultant for ONG RedPapaz (Bogota)
$ws.Range("A2").Value = "External Scientific Consultant"
$ws.Range("B2").Value = "Sep. - Oct. 2022"
$ws.Range("C2").Value = "\href{https://www.redpapaz.org/}{ONG RedPapaz}"
$ws.Range("D2").Value = "Bogotá, Colombia"

# Row 3: International Scientific Consultant for Protect Children (Helsinki)
$ws.Range("A3").Value = "International Scientific Consultant"
$ws.Range("B3").Value = "Jul. - Aug. 2021"
$ws.Range("C3").Value = "\href{https://www.suojellaanlapsia.fi/en}{Protect Children}"
$ws.Range("D3").Value = "Helsink, Finlandia"

# Column E ("why") is no longer used by any data row now - clear its old
# leftover value and shrink the column back down from its old custom width.
$ws.Range("E2").ClearContents()
$ws.Columns("E").ColumnWidth = 8

# Resize columns A and D to fit the new (longer / shorter) content.
$ws.Columns("A").ColumnWidth = 28
$ws.Columns("D").ColumnWidth = 16

# Match the final selection left by the editing session.
$ws.Range("E1").Select() | Out-Null
